$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.001.22"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "'1.822.43"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'309.70"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.4620"
$ws.Range("E7").Value = "  -2.93%  "
$ws.Range("D8").Value = "'0.3638"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").Value = "'0.07285"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").Value = "'0.8647"
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("D11").Value = "'19.84"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").Value = "'1.883.41"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'0.07607"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "'93.24"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'0.000008624"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "'27.438.21"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "'5.164"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'2.118.53"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'151.60"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'1.860"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "'18.16"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").Value = "'2.087"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("D29").Value = "'5.094"
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("D30").Value = "'116.12"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'0.08907"
$ws.Range("D32").Value = "'2.951"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.141"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7270"
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("D35").Value = "'4.426"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("D37").Value = "'2.489"
$ws.Range("E37").Value = "  +4.68%  "
$ws.Range("D38").Value = "'0.05272"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").Value = "'1.072"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").Value = "'0.01918"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").Value = "'2.927"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").Value = "'7.169"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'0.5205"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "'8.264"
$ws.Range("E45").Value = "  -3.50%  "
$ws.Range("D46").Value = "'0.4861"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D49").Value = "'103.15"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").Value = "'1.631"
$ws.Range("E50").Value = "  -3.27%  "
$ws.Range("D51").Value = "'0.06221"
$ws.Range("E51").Value = "  -1.68%  "
